# Update "想去人数" (interest count) figures in column F on the
# "展览" (sheetId 1) and "全部类型" (sheetId 4) worksheets.
# Sheets "演出" and "本地生活" are untouched by this revision.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws4 = $wb.Worksheets.Item("全部类型")

# --- Sheet "展览" ---
$ws1.Range("F2").Value  = 355
$ws1.Range("F4").Value  = 1301
$ws1.Range("F9").Value  = 151
$ws1.Range("F10").Value = 3550
$ws1.Range("F11").Value = 142
$ws1.Range("F14").Value = 52
$ws1.Range("F17").Value = 111
$ws1.Range("F18").Value = 777
$ws1.Range("F19").Value = 215
$ws1.Range("F20").Value = 132
$ws1.Range("F24").Value = 2769
$ws1.Range("F25").Value = 5250
$ws1.Range("F29").Value = 3095
$ws1.Range("F30").Value = 297
$ws1.Range("F31").Value = 2276
$ws1.Range("F34").Value = 88
$ws1.Range("F35").Value = 137
$ws1.Range("F36").Value = 186
$ws1.Range("F38").Value = 40
$ws1.Range("F39").Value = 470
$ws1.Range("F40").Value = 814
$ws1.Range("F42").Value = 8
$ws1.Range("F45").Value = 496

# --- Sheet "全部类型" ---
$ws4.Range("F2").Value  = 355
$ws4.Range("F4").Value  = 1301
$ws4.Range("F9").Value  = 151
$ws4.Range("F10").Value = 3550
$ws4.Range("F11").Value = 142
$ws4.Range("F15").Value = 52
$ws4.Range("F18").Value = 111
$ws4.Range("F19").Value = 777
$ws4.Range("F20").Value = 215
$ws4.Range("F21").Value = 132
$ws4.Range("F25").Value = 2769
$ws4.Range("F26").Value = 5250
$ws4.Range("F30").Value = 3095
$ws4.Range("F31").Value = 297
$ws4.Range("F32").Value = 2276
$ws4.Range("F35").Value = 88
$ws4.Range("F36").Value = 137
$ws4.Range("F37").Value = 186
$ws4.Range("F39").Value = 40
$ws4.Range("F40").Value = 470
$ws4.Range("F41").Value = 814
$ws4.Range("F43").Value = 8
$ws4.Range("F46").Value = 496
